# Auto-committed on 2022/12/23 週五 17:18:57.23
# Adds a new "distItemFirst" lookup row to the DBS sheet, mirroring the
# existing deptCodeFirst / distCodeFirst rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

$ws.Range("A13").Value = "distItemFirst"
$ws.Range("B13").Value = "DistItem %"
$ws.Range("C13").Value = "UnitCode Asc"

$ws.Range("A13").Select() | Out-Null
$ws.Activate() | Out-Null
